$p = $ppt.ActivePresentation
$s = $p.Slides.Item(39)
$sh = $s.Shapes.Item(1)
$tbl = $sh.Table
$cell = $tbl.Cell(1,1)
$cell.Shape.TextFrame.TextRange.Text = "Sources of Spatially Referenced Data"
